# TeamContributions.xlsx - "Add files via upload" re-upload.
#
# The uploaded workbook appends each member's student ID to their name in
# column B, and updates each member's contribution SCORE in column D
# (the three members with a recorded task move from 10 to 33; the member
# with no task recorded - Dale Follows - stays at 0).
#
# Cells are written in the same order the original author must have typed
# them (Dan, then Joseph, then Callum, then Dale) so that new shared-string
# table entries land in the same order as the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Dan Ferguson
$ws.Range("B4").Value = "Dan Ferguson - 40534169"
$ws.Range("D4").Value = 33

# Row 3 - Joseph Fanning
$ws.Range("B3").Value = "Joseph Fanning - 40593072"
$ws.Range("D3").Value = 33

# Row 5 - Callum Hamilton
$ws.Range("B5").Value = "Callum Hamilton - 40591758"
$ws.Range("D5").Value = 33

# Row 6 - Dale Follows (score stays 0, no task logged)
$ws.Range("B6").Value = "Dale Follows - 40606982"

# Leave the saved selection on D3, matching the re-uploaded file.
$ws.Range("D3").Select() | Out-Null
